$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1) Insert a new "Meta description" paragraph right after the Heading1 title ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(2)
$metaXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Down the Rails, a London subway-themed slot game by Pragmatic Play. Play for free and enjoy high-volatility gameplay mechanics.</w:t></w:r></w:p>"
$p2.Range.InsertXML($metaXml)

# --- 2) Remove the old bold title paragraph that used to sit at the bottom ---
$n = $d.Paragraphs.Count
$pOldTitle = $d.Paragraphs.Item($n - 1)
$pOldTitle.Range.Delete()

# --- 3) Replace the trailing italic paragraph's text with the DALLE image prompt ---
$n2 = $d.Paragraphs.Count
$pDesc = $d.Paragraphs.Item($n2)
$promptText = 'Prompt for DALLE: Create a feature image for "Down the Rails" that captures the game''s theme of the London subway and features a happy Maya warrior with glasses. The image should be in a cartoon style and include elements from the game such as the London subway, iconic characters like Shakespeare or Winston Churchill, and bonus features like the End of the Line Bonus game. The Maya warrior should be prominently displayed, perhaps riding the subway or standing in front of Buckingham Palace. The image should be eye-catching and colorful, with a sense of fun and excitement to match the game.'
$descXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$promptText</w:t></w:r></w:p>"
$pDesc.Range.InsertXML($descXml)
